$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.256894556890018
$ws.Range("C2").Value = 0.1859536292247697
$ws.Range("E2").Value = 0.07533167458814383
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.7621162131235195
$ws.Range("H2").Value = 0.8152598442540153
$ws.Range("I2").Value = 0.7262029474934337
$ws.Range("L2").Value = 0.2344177996685488
$ws.Range("N2").Value = 1.162311826403386
$ws.Range("B3").Value = 1.148341961563858
$ws.Range("C3").Value = 0.1670936965257397
$ws.Range("E3").Value = 0.07570507704978624
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.7535796821256326
$ws.Range("H3").Value = 0.8170517083399176
$ws.Range("I3").Value = 0.7307586488425315
$ws.Range("L3").Value = 0.2246696686353289
$ws.Range("N3").Value = 1.176602397373351
$ws.Range("B4").Value = 1.082041398387048
$ws.Range("C4").Value = 0.1554368486202975
$ws.Range("E4").Value = 0.07596398982022112
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.7490281239708025
$ws.Range("H4").Value = 0.8186766823176441
$ws.Range("I4").Value = 0.7340846179621394
$ws.Range("L4").Value = 0.2188154403553
$ws.Range("N4").Value = 1.185904803153992
$ws.Range("B5").Value = 1.055112359836073
$ws.Range("C5").Value = 0.1506672408897316
$ws.Range("E5").Value = 0.07607696025027177
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.7473460131434564
$ws.Range("H5").Value = 0.8194705806767075
$ws.Range("I5").Value = 0.7355725859923936
$ws.Range("L5").Value = 0.2164627653008466
$ws.Range("N5").Value = 1.189828336154168
$ws.Range("B6").Value = 1.050646207812974
$ws.Range("C6").Value = 0.1498740814282655
$ws.Range("E6").Value = 0.07609616977819211
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.7470771064800061
$ws.Range("H6").Value = 0.8196103545509317
$ws.Range("I6").Value = 0.7358276609709868
$ws.Range("L6").Value = 0.2160740964376515
$ws.Range("N6").Value = 1.190487851522942
$ws.Range("B7").Value = 1.08167786278193
$ws.Range("C7").Value = 0.1553726024128537
$ws.Range("E7").Value = 0.07596548315778229
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.7490047402929463
$ws.Range("H7").Value = 0.8186868561586635
$ws.Range("I7").Value = 0.7341041487589308
$ws.Range("L7").Value = 0.218783577874035
$ws.Range("N7").Value = 1.185957179922767
$ws.Range("B8").Value = 1.219393103738696
$ws.Range("C8").Value = 0.1794666275522161
$ws.Range("E8").Value = 0.07545427602702404
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.759029142265419
$ws.Range("H8").Value = 0.815768613269654
$ws.Range("I8").Value = 0.7276637993815385
$ws.Range("L8").Value = 0.2310294101455383
$ws.Range("N8").Value = 1.167129625356107
$ws.Range("B9").Value = 1.492223789773334
$ws.Range("C9").Value = 0.2261098953096621
$ws.Range("E9").Value = 0.07468669110313897
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.7841992762957659
$ws.Range("H9").Value = 0.8142218554891087
$ws.Range("I9").Value = 0.7192460246762238
$ws.Range("L9").Value = 0.256086290170316
$ws.Range("N9").Value = 1.134398828888727
$ws.Range("B10").Value = 1.694364019382704
$ws.Range("C10").Value = 0.2600187314510265
$ws.Range("E10").Value = 0.07426557213452512
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.806108420687778
$ws.Range("H10").Value = 0.8156494614578236
$ws.Range("I10").Value = 0.7156524840102065
$ws.Range("L10").Value = 0.27513654784984
$ws.Range("N10").Value = 1.112907004904017
$ws.Range("B11").Value = 1.786692085938341
$ws.Range("C11").Value = 0.2753686875401797
$ws.Range("E11").Value = 0.07410493376704785
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.8168296879691184
$ws.Range("H11").Value = 0.8168597328162548
$ws.Range("I11").Value = 0.714585427558994
$ws.Range("L11").Value = 0.2839434634916103
$ws.Range("N11").Value = 1.103685076592527
$ws.Range("B12").Value = 1.821707712018508
$ws.Range("C12").Value = 0.2811705825341164
$ws.Range("E12").Value = 0.07404854565874963
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.8209990149046291
$ws.Range("H12").Value = 0.8173989912404807
$ws.Range("I12").Value = 0.7142633999801546
$ws.Range("L12").Value = 0.2872987248687053
$ws.Range("N12").Value = 1.100272817201557
$ws.Range("B13").Value = 1.814164124865727
$ws.Range("C13").Value = 0.2799215204250913
$ws.Range("E13").Value = 0.07406049236991841
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.820096195208265
$ws.Range("H13").Value = 0.8172792468982948
$ws.Range("I13").Value = 0.714329098852609
$ws.Range("L13").Value = 0.286575207322457
$ws.Range("N13").Value = 1.101004153702661
$ws.Range("B14").Value = 1.789571784592624
$ws.Range("C14").Value = 0.2758462291862145
$ws.Range("E14").Value = 0.07410020568815412
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.8171705030274552
$ws.Range("H14").Value = 0.8169024737127586
$ws.Range("I14").Value = 0.7145572875013002
$ws.Range("L14").Value = 0.2842190964432234
$ws.Range("N14").Value = 1.103402746334588
$ws.Range("B15").Value = 1.774515146128351
$ws.Range("C15").Value = 0.2733485901559618
$ws.Range("E15").Value = 0.07412510957568408
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.8153927077897265
$ws.Range("H15").Value = 0.8166822411001817
$ws.Range("I15").Value = 0.7147077564585445
$ws.Range("L15").Value = 0.2827785512775307
$ws.Range("N15").Value = 1.104882360109244
$ws.Range("B16").Value = 1.688337616969818
$ws.Range("C16").Value = 0.259014062496334
$ws.Range("E16").Value = 0.07427669215076094
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.8054230231206958
$ws.Range("H16").Value = 0.8155816791319239
$ws.Range("I16").Value = 0.7157336768583846
$ws.Range("L16").Value = 0.2745638311671428
$ws.Range("N16").Value = 1.113520852530201
$ws.Range("B17").Value = 1.635565636474553
$ws.Range("C17").Value = 0.2502010205912484
$ws.Range("E17").Value = 0.07437760121353598
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.7995008721017172
$ws.Range("H17").Value = 0.815050383075544
$ws.Range("I17").Value = 0.7165087320657051
$ws.Range("L17").Value = 0.2695604562684224
$ws.Range("N17").Value = 1.118962459873757
$ws.Range("B18").Value = 1.60524778299316
$ws.Range("C18").Value = 0.2451248981806202
$ws.Range("E18").Value = 0.07443855312061309
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.796165588668444
$ws.Range("H18").Value = 0.8147975666764609
$ws.Range("I18").Value = 0.7170079295618521
$ws.Range("L18").Value = 0.2666959044121455
$ws.Range("N18").Value = 1.122144543247842
$ws.Range("B19").Value = 1.594988745212788
$ws.Range("C19").Value = 0.2434049909644784
$ws.Range("E19").Value = 0.07445969064679048
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.7950484829585349
$ws.Range("H19").Value = 0.8147210215910974
$ws.Range("I19").Value = 0.7171861092469598
$ws.Range("L19").Value = 0.2657282927410449
$ws.Range("N19").Value = 1.123230906885205
$ws.Range("B20").Value = 1.641179666642529
$ws.Range("C20").Value = 0.2511399178900149
$ws.Range("E20").Value = 0.07436655796466418
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.8001239424429514
$ws.Range("H20").Value = 0.8151014766298772
$ws.Range("I20").Value = 0.716420695771923
$ws.Range("L20").Value = 0.2700917018848088
$ws.Range("N20").Value = 1.11837778671368
$ws.Range("B21").Value = 1.796793718518643
$ws.Range("C21").Value = 0.2770435338468928
$ws.Range("E21").Value = 0.07408842040972452
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.8180268730999387
$ws.Range("H21").Value = 0.8170109415860054
$ws.Range("I21").Value = 0.7144880330135877
$ws.Range("L21").Value = 0.2849105929511211
$ws.Range("N21").Value = 1.102696052630563
$ws.Range("B22").Value = 1.898805046760003
$ws.Range("C22").Value = 0.2939101090444467
$ws.Range("E22").Value = 0.0739325311560517
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.8303656482892166
$ws.Range("H22").Value = 0.818730898365061
$ws.Range("I22").Value = 0.7137033324797954
$ws.Range("L22").Value = 0.2947137740146246
$ws.Range("N22").Value = 1.092912814599323
$ws.Range("B23").Value = 1.844331681258893
$ws.Range("C23").Value = 0.2849138473928008
$ws.Range("E23").Value = 0.07401336511859036
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.8237215229032699
$ws.Range("H23").Value = 0.8177696342916363
$ws.Range("I23").Value = 0.7140782302947812
$ws.Range("L23").Value = 0.2894708120406193
$ws.Range("N23").Value = 1.098091667115064
$ws.Range("B24").Value = 1.638641496702064
$ws.Range("C24").Value = 0.2507154716693378
$ws.Range("E24").Value = 0.0743715414659345
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.7998420361046783
$ws.Range("H24").Value = 0.8150782132999694
$ws.Range("I24").Value = 0.7164603300195296
$ws.Range("L24").Value = 0.2698514885070722
$ws.Range("N24").Value = 1.118641950390007
$ws.Range("B25").Value = 1.418118991434596
$ws.Range("C25").Value = 0.2135555167519954
$ws.Range("E25").Value = 0.07486923506585264
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.7767941610114519
$ws.Range("H25").Value = 0.8194705806767075
$ws.Range("I25").Value = 0.7147077564585445
$ws.Range("L25").Value = 0.2491955609159362
$ws.Range("N25").Value = 1.142804713412488
